$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 92 (shifts existing rows 92-104 down to 93-105),
# copying formatting from the row above as Excel normally does.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly price point.
$ws.Cells.Item(92, 1).Value = 7
$ws.Cells.Item(92, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(92, 3).Value = "Ñuble"
$ws.Cells.Item(92, 4).Value = 45132
$ws.Cells.Item(92, 5).Value = 16
$ws.Cells.Item(92, 6).Value = 100112044
$ws.Cells.Item(92, 7).Value = "Perejil"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 250
$ws.Cells.Item(92, 11).Value = 1500
$ws.Cells.Item(92, 12).Value = 1500
$ws.Cells.Item(92, 13).Value = 1500
$ws.Cells.Item(92, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(92, 15).Value = "Región de Ñuble"
$ws.Cells.Item(92, 16).Value = 1500
$ws.Cells.Item(92, 17).Value = 1
$ws.Cells.Item(92, 18).Value = "Hortaliza"
